$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new time entry for row 26 (date already present in A26)
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = "Update UI and animation of the assistant view."

# Copy style (font size 20) from row 25 B/C cells to the new B26/C26 cells
$ws.Range("B25").Copy()
$ws.Range("B26").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C25").Copy()
$ws.Range("C26").PasteSpecial(-4122) # xlPasteFormats

# Update the SUM formula to include the new row
$ws.Range("B42").Formula = "=SUM(B2:B26)"

# Update the selected cell in the view
$ws.Range("C26").Select()
